$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.146.42'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '3.180.98'
$ws.Range("E3").Value = '  -3.96%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.60%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.178.46'
$ws.Range("E8").Value = '  -3.99%  '
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("E10").Value = '  -6.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.22'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.453'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.44%  '
$ws.Range("E13").Value = '  -4.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").Value = '3.704.96'
$ws.Range("E15").Value = '  -3.95%  '
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '3.185.58'
$ws.Range("D18").Value = '63.062.17'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.56'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.698'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  -3.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.102'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").Value = '0.0₃0708'
$ws.Range("E38").Value = '  -5.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0388'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '406.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.19%  '
$ws.Range("E42").Value = '  -3.64%  '
$ws.Range("E43").Value = '  -5.97%  '
$ws.Range("D44").Value = '2.815.30'
$ws.Range("E44").Value = '  -9.61%  '
$ws.Range("E45").Value = '  -5.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.58%  '
$ws.Range("E51").Value = '  -2.06%  '
